# LOBSTAHS_lipid_class_rt_windows.xlsx - "Update files with title of paper as submitted"
#
# On the "Notes" sheet:
#   1. The whole notes table (previously parked in columns E:G, presumably so
#      it wouldn't print/overlap column D) is moved back to columns A:C.
#   2. A new line citing the published paper is inserted right after the
#      "Latest versions of all scripts..." line, pushing the remaining notes
#      (including the file-history table) down by one row.
#   3. Window/selection cosmetics are updated to match.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Notes")

# --- 1. Shift the notes table from E:G back to A:C ------------------------
$src = $ws2.Range("E1:G17")
$dst = $ws2.Range("A1")
$src.Copy($dst)
$src.Clear()

# Preserve the (now orphaned) source columns' widths on their new home
$ws2.Columns.Item(1).ColumnWidth = 8.498697916666666   # -> stored width 9.33203125
$ws2.Columns.Item(2).ColumnWidth = 31.830729166666668  # -> stored width 32.6640625

# --- 2. Insert the new citation line at row 7 ------------------------------
$ws2.Rows.Item(7).Insert()
$ws2.Range("A7").Value = "See Collins, J.R., B.R. Edwards, H.F. Fredricks, and B.A.S. Van Mooy, 2016, ""LOBSTAHS: A Novel Lipidomics Strategy for Semi-Untargeted Discovery and Identification of Oxidative Stress Biomarkers"""

# --- 3. Selection / view cosmetics -----------------------------------------
$ws2.Range("H18").Select()
